$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Explicit target values for columns D (x_corrSteps), F (x_nrSteps) and H (alienID)
# for the affected trial rows, as described by the diff.
$updates = @{
    2  = @{ D = 4; F = -3; H = 46 }
    9  = @{ D = 2; F = -3; H = 46 }
    12 = @{ D = 4; F = -3; H = 46 }
    20 = @{ D = 2; F = -3; H = 46 }
    22 = @{ D = 6; F = -3; H = 46 }
    28 = @{ D = 2; F = -3; H = 46 }
}

foreach ($r in $updates.Keys) {
    $vals = $updates[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("F$r").Value = $vals.F
    $ws.Range("H$r").Value = $vals.H
}

# Update the active selection to D28, matching the saved sheet view.
$ws.Range("D28").Select()
